$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '54.995.43'
$ws.Range("E2").Value = '  +1.16%  '

$ws.Range("D3").Value = '2.290.22'
$ws.Range("E3").Value = '  +0.30%  '

$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '507.15'
$ws.Range("E5").Value = '  +0.81%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.47'
$ws.Range("E6").Value = '  -0.06%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.997'
$ws.Range("E7").Value = '  -0.13%  '

$ws.Range("E8").Value = '  +0.19%  '

$ws.Range("D9").Value = '2.312.40'
$ws.Range("E9").Value = '  +0.80%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0972'
$ws.Range("E10").Value = '  +1.64%  '

$ws.Range("E11").Value = '  +1.81%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.342'
$ws.Range("E12").Value = '  +2.52%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.94'
$ws.Range("E13").Value = '  +4.29%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.66'
$ws.Range("E14").Value = '  +2.67%  '

$ws.Range("D15").Value = '2.697.94'
$ws.Range("E15").Value = '  +0.19%  '

$ws.Range("D16").Value = '54.725.63'
$ws.Range("E16").Value = '  +0.74%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000132'
$ws.Range("E17").Value = '  +1.43%  '

$ws.Range("D18").Value = '2.268.18'
$ws.Range("E18").Value = '  -1.52%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.61'
$ws.Range("E19").Value = '  +3.33%  '

$ws.Range("E20").Value = '  +1.53%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.68'
$ws.Range("E21").Value = '  +3.98%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '308.41'
$ws.Range("E22").Value = '  +1.07%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.998'
$ws.Range("E23").Value = '  -0.13%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '60.45'
$ws.Range("E24").Value = '  -2.67%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.01'
$ws.Range("E25").Value = '  +1.05%  '

$ws.Range("E26").Value = '  -0.06%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.51'
$ws.Range("E27").Value = '  +2.16%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '171.72'
$ws.Range("E28").Value = '  -1.54%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '6.13'
$ws.Range("E29").Value = '  +1.98%  '

$ws.Range("B30").Value = 'Fetch.AI'
$ws.Range("C30").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.15'
$ws.Range("E30").Value = '  +6.61%  '

$ws.Range("B31").Value = 'PEPE'
$ws.Range("C31").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D31").Value = '0.0₃0707'
$ws.Range("E31").Value = '  +2.08%  '

$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.63'
$ws.Range("E32").Value = '  +0.38%  '

$ws.Range("E33").Value = '  -0.01%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '18.04'
$ws.Range("E34").Value = '  +1.45%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.996'
$ws.Range("E35").Value = '  -0.06%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.910'
$ws.Range("E36").Value = '  -3.40%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.22'
$ws.Range("E37").Value = '  +1.01%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.84'
$ws.Range("E38").Value = '  +1.63%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.64'
$ws.Range("E39").Value = '  +1.58%  '

$ws.Range("E40").Value = '  +0.58%  '

$ws.Range("E41").Value = '  +1.51%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '133.50'
$ws.Range("E42").Value = '  +6.51%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.43'
$ws.Range("E43").Value = '  +0.88%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.88'
$ws.Range("E44").Value = '  +1.20%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '253.91'
$ws.Range("E45").Value = '  +5.31%  '

$ws.Range("E46").Value = '  +1.16%  '

$ws.Range("E47").Value = '  +1.89%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.554'
$ws.Range("E48").Value = '  +0.94%  '

$ws.Range("E49").Value = '  +0.89%  '

$ws.Range("E50").Value = '  +0.86%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '10.82'
$ws.Range("E51").Value = '  +0.47%  '
